$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 31
$ws.Range("B10").Value = "Update index2.py"
$ws.Range("C10").Value = "riya-morankar"
$ws.Range("D10").Value = "N/A"
$ws.Range("E10").Value = "edit2 to main"

# Force the date-looking text to stay as literal text instead of being
# auto-converted into a date serial number by Excel's input parser.
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "2025-06-18"
$ws.Range("F10").Style = "Normal"
